$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "22.70" or "0.996"
# keep their exact string representation instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "61.590.28"
$ws.Range("E2").Value = "  +12.06%  "

# Row 3
$ws.Range("D3").Value = "2.650.31"
$ws.Range("E3").Value = "  +13.39%  "

# Row 4
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.52%  "

# Row 5
$ws.Range("D5").Value = "517.48"
$ws.Range("E5").Value = "  +9.16%  "

# Row 6
$ws.Range("D6").Value = "162.11"
$ws.Range("E6").Value = "  +12.07%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  +2.12%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "0.995"
$ws.Range("E8").Value = "  -0.41%  "

# Row 9
$ws.Range("D9").Value = "2.683.86"
$ws.Range("E9").Value = "  +14.47%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  +12.08%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  +13.47%  "

# Row 12
$ws.Range("E12").Value = "  +8.94%  "

# Row 13
$ws.Range("E13").Value = "  +1.84%  "

# Row 14
$ws.Range("D14").Value = "3.137.95"
$ws.Range("E14").Value = "  +14.11%  "

# Row 15
$ws.Range("D15").Value = "61.260.01"
$ws.Range("E15").Value = "  +11.28%  "

# Row 16
$ws.Range("D16").Value = "22.70"
$ws.Range("E16").Value = "  +14.58%  "

# Row 17
$ws.Range("E17").Value = "  +11.71%  "

# Row 18
$ws.Range("D18").Value = "2.693.46"
$ws.Range("E18").Value = "  +14.71%  "

# Row 19
$ws.Range("D19").Value = "4.86"
$ws.Range("E19").Value = "  +6.58%  "

# Row 20
$ws.Range("D20").Value = "356.69"
$ws.Range("E20").Value = "  +14.02%  "

# Row 21
$ws.Range("D21").Value = "10.69"
$ws.Range("E21").Value = "  +12.74%  "

# Row 22
$ws.Range("E22").Value = "  +11.29%  "

# Row 23
$ws.Range("E23").Value = "  +0.10%  "

# Row 24
$ws.Range("D24").Value = "61.15"
$ws.Range("E24").Value = "  +8.83%  "

# Row 25
$ws.Range("D25").Value = "0.432"
$ws.Range("E25").Value = "  +10.26%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.801.37"
$ws.Range("E26").Value = "  +14.46%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.170"
$ws.Range("E27").Value = "  +12.08%  "

# Row 28
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +1.01%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0872"
$ws.Range("E29").Value = "  +18.72%  "

# Row 30
$ws.Range("D30").Value = "7.66"
$ws.Range("E30").Value = "  +9.20%  "

# Row 31
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("E32").Value = "  +9.94%  "

# Row 33
$ws.Range("D33").Value = "158.99"
$ws.Range("E33").Value = "  +9.41%  "

# Row 34
$ws.Range("E34").Value = "  +9.69%  "

# Row 35
$ws.Range("D35").Value = "5.73"
$ws.Range("E35").Value = "  +13.05%  "

# Row 36
$ws.Range("E36").Value = "  +13.74%  "

# Row 37
$ws.Range("E37").Value = "  +13.32%  "

# Row 38
$ws.Range("D38").Value = "0.884"
$ws.Range("E38").Value = "  +10.28%  "

# Row 39
$ws.Range("E39").Value = "  +15.55%  "

# Row 40
$ws.Range("E40").Value = "  +37.64%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").Value = "  +13.76%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "305.59"
$ws.Range("E42").Value = "  +23.45%  "

# Row 43
$ws.Range("D43").Value = "36.10"
$ws.Range("E43").Value = "  +7.66%  "

# Row 44
$ws.Range("D44").Value = "0.647"
$ws.Range("E44").Value = "  +12.42%  "

# Row 45
$ws.Range("E45").Value = "  +14.88%  "

# Row 46
$ws.Range("D46").Value = "0.103"
$ws.Range("E46").Value = "  +1.56%  "

# Row 47
$ws.Range("D47").Value = "20.42"
$ws.Range("E47").Value = "  +23.42%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "5.06"
$ws.Range("E48").Value = "  +16.28%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "0.993"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("D50").Value = "0.0241"
$ws.Range("E50").Value = "  +9.94%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.92"
$ws.Range("E51").Value = "  +20.11%  "

# Restore default (Normal) style for column D so no stray number-format style remains on cells
$ws.Range("D2:D51").Style = "Normal"